# Add 5 new medicine rows (22-26) to the ListaCuMedicamente sheet, using the
# same look & feel as the existing rows, then move the active selection to G27
# (the first empty row after the newly-added data), matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Copy the per-column formatting (fill/border/alignment) from the nearest
#    existing "complete" rows so the new rows look consistent with the rest
#    of the table. Row 6 and Row 13 alternate the same striped style used
#    throughout the sheet; Row 4 supplies the wrap-text variant needed for
#    column D on row 26.
# ---------------------------------------------------------------------------
$ws.Range("A6:H6").Copy($ws.Range("A22:H22"))
$ws.Range("A13:H13").Copy($ws.Range("A23:H23"))
$ws.Range("A6:H6").Copy($ws.Range("A24:H24"))
$ws.Range("A13:H13").Copy($ws.Range("A25:H25"))
$ws.Range("A4:H4").Copy($ws.Range("A26:H26"))
# Column F on row 26 holds a plain number (like row 6's F), not a shared
# string placeholder, so pull that specific cell format from row 6 instead.
$ws.Cells.Item(6, 6).Copy($ws.Cells.Item(26, 6))

# ---------------------------------------------------------------------------
# 2) Row 22 - Egistrozol 1mg x 30 comprimate filmate
#    (values entered in column order A,B,C,D,E)
# ---------------------------------------------------------------------------
$ws.Cells.Item(22, 1).Value = "Egistrozol 1mg x 30 comprimate filmate"
$ws.Cells.Item(22, 2).Value = "Imagine_21.jpg"
$ws.Cells.Item(22, 3).Value = "Doza recomandată este de un comprimat administrat o dată pe z"
$ws.Cells.Item(22, 4).Value = "Antineoplazice"
$ws.Cells.Item(22, 5).Value = "Egistrozol 1 mg comprimate filmate conţine o substanţă numită anastrozol. Acesta aparţine unui grup de medicamente, numite inhibitori de aromatază. Egistrozol este utilizat în tratamentul cancerului de sân la femeile care au trecut la menopauză."
$ws.Cells.Item(22, 6).Value = 12
$ws.Cells.Item(22, 7).Value = "RON"
$ws.Rows.Item(22).RowHeight = 60

# ---------------------------------------------------------------------------
# 3) Row 23 - Decaris 150mg, 1 comprimat
#    (values entered in the order A,B,E,D,C -- matches shared-string order)
# ---------------------------------------------------------------------------
$ws.Cells.Item(23, 1).Value = "Decaris 150mg, 1 comprimat"
$ws.Cells.Item(23, 2).Value = "Imagine_22.jpg"
$ws.Cells.Item(23, 5).Value = "DECARIS face parte din grupa antinematode, derivaţi de imidazotiazol; are acţiune vermicidă prin paralizia musculaturii viermilor intestinali.`nDECARIS este indicat în infestări cu: Ascaris lumbricoides, Necator americanus, Strongyloides stercocalis, Trichostrongylus colubriformis şi Ankylostoma duodenale"
$ws.Cells.Item(23, 4).Value = "Produs antiparazitar"
$ws.Cells.Item(23, 3).Value = "1 comprimat de 150 mg în doză unică.Este preferabil ca administrarea medicamentului să se facă seara, la culcare. La nevoie, tratamentul se va repeta după o pauză de 14 zile."
$ws.Cells.Item(23, 6).Value = 21
$ws.Cells.Item(23, 7).Value = "RON"
$ws.Rows.Item(23).RowHeight = 90

# ---------------------------------------------------------------------------
# 4) Row 24 - Zentel suspensie 200mg/5ml, 10ml
#    (values entered in order A,B,C,D,E -- D reuses the existing
#    "Produs antiparazitar" string created for row 23)
# ---------------------------------------------------------------------------
$ws.Cells.Item(24, 1).Value = "Zentel suspensie 200mg/5ml, 10ml"
$ws.Cells.Item(24, 2).Value = "Imagine_23.jpg"
$ws.Cells.Item(24, 3).Value = "Medicul dumneavoastră vă va recomanda ce doză de suspensie ZENTEL este necesară zilnic şi cât timp trebuie să luaţi ZENTEL. Doza depinde de greutatea sau vârsta dumneavoastră şi de tipul şi severitatea infecţiei."
$ws.Cells.Item(24, 4).Value = "Produs antiparazitar"
$ws.Cells.Item(24, 5).Value = "ZENTEL este o suspensie care trebuie administrată pe cale orală.`nZENTEL este un carbamat de benzimidazol, compus care aparţine grupei de medicamente antihelmintice şi antiparazitare.`nZENTEL este utilizat pentru tratamentul unei varietăţi largi de afecţiuni intestinale determinate de viermi şi paraziţi.`nSe presupune că ZENTEL elimină viermii sau paraziţii prin afectarea metabolismului acestora ceea ce determina moartea lor. Sunt afectate atât ouăle şi larvele cât şi paraziţii adulţi."
$ws.Cells.Item(24, 6).Value = 45
$ws.Cells.Item(24, 7).Value = "RON"
$ws.Rows.Item(24).RowHeight = 120

# ---------------------------------------------------------------------------
# 5) Row 25 - Arnigel, 45 g
#    (values entered in order A,B,C,D,E)
# ---------------------------------------------------------------------------
$ws.Cells.Item(25, 1).Value = "Arnigel, 45 g"
$ws.Cells.Item(25, 2).Value = "Imagine_24.jpg"
$ws.Cells.Item(25, 3).Value = "De 2-3 ori pe zi, la adulti si la copii mai mari de 1 an."
$ws.Cells.Item(25, 4).Value = "Unguent"
$ws.Cells.Item(25, 5).Value = "Arnigel este un medicament homeopat utilizat in mod traditional in tratamentul local in traumatologia benigna in absenta plagilor (echimoze, contuzii, oboseala musculara) la adulti si copii cu varsta peste 1 an.`n"
$ws.Cells.Item(25, 6).Value = 14
$ws.Cells.Item(25, 7).Value = "RON"
$ws.Rows.Item(25).RowHeight = 60

# ---------------------------------------------------------------------------
# 6) Row 26 - Sulfat de bariu suspensie orala, 95g MED
#    (values entered in the order A,B,D,C,E -- matches shared-string order)
# ---------------------------------------------------------------------------
$ws.Cells.Item(26, 1).Value = "Sulfat de bariu suspensie orala, 95g MED"
$ws.Cells.Item(26, 2).Value = "Imagine_25.jpg"
$ws.Cells.Item(26, 4).Value = "Pulbere pentru suspensie orala"
$ws.Cells.Item(26, 3).Value = "Continutul unui flacon se amesteca cu apa, se agita bine si se bea pe stomacul gol, dupa indicatiile medicului radioimagist."
$ws.Cells.Item(26, 5).Value = "SULFAT DE BARIU apartine unei clase de medicamente cunoscute sub numele de medii de contrast pentru radiologie, fara iod. Este folosit pentru explorarea tubului digestiv."
$ws.Cells.Item(26, 6).Value = 22
$ws.Cells.Item(26, 7).Value = "RON"
$ws.Rows.Item(26).RowHeight = 60

# ---------------------------------------------------------------------------
# 7) Move the active selection to G27, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("G27").Select()
